$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 0
    20 = 3
    21 = 1
    22 = 8
    23 = 4
    24 = 5
    25 = 4
    26 = 5
    27 = 2
    28 = 2
    29 = 4
    30 = 1
    31 = 4
    32 = 2
    33 = 5
    34 = 1
    35 = 0
    36 = 1
    37 = 2
    38 = 0
    39 = 0
    40 = 1
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 3
    48 = 1
    49 = 0
    50 = 1
    51 = 2
    52 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
